$wb = $excel.ActiveWorkbook

# --- Content edit: translation keys renamed -------------------------------
# "survey" sheet header row: display.text -> display.prompt.text,
#                             display.hint -> display.hint.text
$survey = $wb.Worksheets.Item("survey")
$survey.Range("D1").Value = "display.prompt.text"
$survey.Range("E1").Value = "display.hint.text"

# "settings" sheet: the display.title setting is renamed display.title.text
$settings = $wb.Worksheets.Item("settings")
$settings.Range("C1").Value = "display.title.text"

# --- View / selection state ------------------------------------------------
# Park the cursor on every other sheet first (selecting a range on a sheet
# implicitly activates it), then activate + select "survey" last so it ends
# up as the active tab with the right cell selected.
$model = $wb.Worksheets.Item("model")
[void]$model.Range("B7").Select()

$choices = $wb.Worksheets.Item("choices")
[void]$choices.Range("A12").Select()

[void]$settings.Range("C2").Select()

[void]$survey.Activate()
[void]$survey.Range("E2").Select()
